$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 319, shifting existing rows 319:346 down to 320:347
$ws.Rows(319).Insert()

# Populate the newly inserted row 319 with the new weekly data point
$ws.Cells.Item(319, 1).Value = 7
$ws.Cells.Item(319, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(319, 3).Value = "Ñuble"
$ws.Cells.Item(319, 4).Value = 45166
$ws.Cells.Item(319, 5).Value = 16
$ws.Cells.Item(319, 6).Value = 100112032
$ws.Cells.Item(319, 7).Value = "Zapallo italiano"
$ws.Cells.Item(319, 8).Value = "Sin especificar"
$ws.Cells.Item(319, 9).Value = "Primera"
$ws.Cells.Item(319, 10).Value = 100
$ws.Cells.Item(319, 11).Value = 14000
$ws.Cells.Item(319, 12).Value = 14000
$ws.Cells.Item(319, 13).Value = 14000
$ws.Cells.Item(319, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(319, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(319, 16).Value = 280
$ws.Cells.Item(319, 17).Value = 50
$ws.Cells.Item(319, 18).Value = "Hortaliza"
